$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update the "AMS [N]" measurements in column C (formulas in column D / the
# linked chart will recalculate automatically from these new inputs).
$ws.Range("C2").Value = 99.737080000000006
$ws.Range("C3").Value = 282.56619999999998
$ws.Range("C4").Value = 570.31399999999996
$ws.Range("C5").Value = 663.96690000000001
$ws.Range("C6").Value = 1392.172
$ws.Range("C7").Value = 1388.2860000000001
$ws.Range("C8").Value = 2892.1970000000001
$ws.Range("C9").Value = 2474.6509999999998

# Update the saved selection on the sheet.
$ws.Range("G14").Select()

$wb.Save()
